$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 6
$ws.Range("Q2").Value = 3.6
$ws.Range("R2").Value = 6.574989197057679

# Row 3
$ws.Range("L3").Value = 8
$ws.Range("Q3").Value = 4.8
$ws.Range("R3").Value = 7.8

# Row 4
$ws.Range("L4").Value = 9
$ws.Range("Q4").Value = 5.4
$ws.Range("R4").Value = 8.4

# Row 5
$ws.Range("L5").Value = 6
$ws.Range("Q5").Value = 3.6
$ws.Range("R5").Value = 7.35820382759259

# Row 6
$ws.Range("L6").Value = 6
$ws.Range("Q6").Value = 3.6
$ws.Range("R6").Value = 5.480797077977883

# Row 7
$ws.Range("L7").Value = 3
$ws.Range("Q7").Value = 1.8
$ws.Range("R7").Value = 2.8

# Row 8
$ws.Range("L8").Value = 3
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 4.422459331201854

# Row 9
$ws.Range("L9").Value = 1
$ws.Range("Q9").Value = 0.6
$ws.Range("R9").Value = 3.984053487157619

# Row 10
$ws.Range("L10").Value = 7
$ws.Range("Q10").Value = 4.2
$ws.Range("R10").Value = 8.2

# Row 11
$ws.Range("L11").Value = 4
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 5.4

# Row 12
$ws.Range("L12").Value = 7
$ws.Range("Q12").Value = 4.2
$ws.Range("R12").Value = 6.2

# Row 14
$ws.Range("L14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 4

# Row 15
$ws.Range("L15").Value = 2
$ws.Range("Q15").Value = 1.2
$ws.Range("R15").Value = 3.95820382759259

# Row 16
$ws.Range("L16").Value = 10
$ws.Range("Q16").Value = 6
$ws.Range("R16").Value = 7

# Row 17
$ws.Range("L17").Value = 6
$ws.Range("Q17").Value = 3.6
$ws.Range("R17").Value = 7.6

# Row 18
$ws.Range("L18").Value = 1
$ws.Range("Q18").Value = 0.6
$ws.Range("R18").Value = 0.6

# Row 19
$ws.Range("L19").Value = 7
$ws.Range("Q19").Value = 4.2
$ws.Range("R19").Value = 8.2

# Row 20
$ws.Range("L20").Value = 5
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 7

# Row 21
$ws.Range("L21").Value = 6
$ws.Range("Q21").Value = 3.6
$ws.Range("R21").Value = 7.6

# Row 23
$ws.Range("L23").Value = 7
$ws.Range("Q23").Value = 4.2
$ws.Range("R23").Value = 6.822459331201855

# Row 24
$ws.Range("L24").Value = 7
$ws.Range("Q24").Value = 4.2
$ws.Range("R24").Value = 7.95820382759259

# Row 25
$ws.Range("L25").Value = 6
$ws.Range("Q25").Value = 3.6
$ws.Range("R25").Value = 7.6

# Row 26
$ws.Range("L26").Value = 3
$ws.Range("Q26").Value = 1.8
$ws.Range("R26").Value = 5.250446498266093

# Row 27
$ws.Range("L27").Value = 9
$ws.Range("Q27").Value = 5.4
$ws.Range("R27").Value = 9.15820382759259

# Row 28
$ws.Range("L28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = 4

# Row 29
$ws.Range("L29").Value = 1
$ws.Range("Q29").Value = 0.6
$ws.Range("R29").Value = 4.45640103125154

# Row 30
$ws.Range("L30").Value = 10

# Row 31
$ws.Range("L31").Value = 4
$ws.Range("Q31").Value = 2.4
$ws.Range("R31").Value = 6.227987167064239

# Row 32
$ws.Range("L32").Value = 9
$ws.Range("Q32").Value = 5.4
$ws.Range("R32").Value = 7.4

# Row 33
$ws.Range("L33").Value = 10

# Row 34
$ws.Range("Q34").Value = 0.6
$ws.Range("R34").Value = 4.6

# Row 35
$ws.Range("L35").Value = 9
$ws.Range("Q35").Value = 5.4
$ws.Range("R35").Value = 6.4

# Row 36
$ws.Range("L36").Value = 5
$ws.Range("Q36").Value = 3
$ws.Range("R36").Value = 3
